$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

# Columns A:N, P:Q hold text (even digit-only / date-look-alike values must
# stay text, matching the rest of the sheet); column O holds a real number.
$textRng = $ws.Range("A$row`:N$row")
$textRng.NumberFormat = "@"
$textRng2 = $ws.Range("P$row`:Q$row")
$textRng2.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "1/1/2025"
$ws.Cells.Item($row, 2).Value = "5:38:29 pm"
$ws.Cells.Item($row, 3).Value = "010125173829"
$ws.Cells.Item($row, 4).Value = "order_PeAp7aj0j7sjOl"
$ws.Cells.Item($row, 5).Value = "21"
$ws.Cells.Item($row, 6).Value = "21B81A05V9"
$ws.Cells.Item($row, 7).Value = "SAMRATH REDDY"
$ws.Cells.Item($row, 8).Value = "CSE"
$ws.Cells.Item($row, 9).Value = "E"
$ws.Cells.Item($row, 10).Value = "+917981455290"
$ws.Cells.Item($row, 11).Value = "samrathreddy04@gmail.com"
$ws.Cells.Item($row, 12).Value = "CollegeFee"
$ws.Cells.Item($row, 13).Value = "IV"
# Column N ("FeeSem") is blank on this row, but must be a literal empty
# text value (not an absent/empty cell) to match the rest of the sheet -
# a lone apostrophe is Excel's "force text" marker and collapses to "".
$ws.Cells.Item($row, 14).Value = "'"
$ws.Cells.Item($row, 15).Value = 120000
$ws.Cells.Item($row, 16).Value = "wallet"
$ws.Cells.Item($row, 17).Value = "Verification in progress..."

$textRng.ClearFormats()
$textRng2.ClearFormats()
